$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F header (row 1): new timestamp, same style as the other headers ---
$ws.Cells.Item(1, 6).Value = "2026-01-27 19:19:40"
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)  # xlPasteFormats - copy style/format only

# --- Data rows 2-100: column F mirrors column E (numeric prices) ---
for ($r = 2; $r -le 100; $r++) {
    $e = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 6).Value = $e
}

# --- Rows 101-204: column F is an empty cell, same as columns D/E for those rows ---
for ($r = 101; $r -le 204; $r++) {
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 6).Font.Bold = $false
}
